$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1753503333333334
$ws.Range("H2").Value = 0.526051
$ws.Range("I2").Value = 0.002060424806616007
$ws.Range("J2").Value = 0.002147484514575959
$ws.Range("M2").Value = 1.772429333333333
$ws.Range("N2").Value = 5.317288
$ws.Range("O2").Value = 0.2841077240522499
$ws.Range("P2").Value = 0.2841077240522499
$ws.Range("Q2").Value = 0.3107960744097778
$ws.Range("R2").Value = 2.797164669688
$ws.Range("S2").Value = 0.0005853826023884709
$ws.Range("T2").Value = 0.0006101169378736265

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1753503333333334
$ws.Range("H3").Value = 0.526051
$ws.Range("I3").Value = 0.002060424806616007
$ws.Range("J3").Value = 0.002147484514575959
$ws.Range("O3").Value = 0.01265586741773478
$ws.Range("P3").Value = 0.01265586741773478
$ws.Range("Q3").Value = 0.01384472711822222
$ws.Range("R3").Value = 0.124602544064
$ws.Range("S3").Value = [double]"2.6076463176744E-05"
$ws.Range("T3").Value = [double]"2.717827929811187E-05"

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1753503333333334
$ws.Range("H4").Value = 0.526051
$ws.Range("I4").Value = 0.002060424806616007
$ws.Range("J4").Value = 0.002147484514575959
$ws.Range("M4").Value = 2.471359333333333
$ws.Range("N4").Value = 7.414078
$ws.Range("O4").Value = 0.3961411957610453
$ws.Range("P4").Value = 0.3961411957610453
$ws.Range("Q4").Value = 0.4333536828864445
$ws.Range("R4").Value = 3.900183145978001
$ws.Range("S4").Value = 0.0008162191466685855
$ws.Range("T4").Value = 0.0008507070834824484

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1753503333333334
$ws.Range("H5").Value = 0.526051
$ws.Range("I5").Value = 0.002060424806616007
$ws.Range("J5").Value = 0.002147484514575959
$ws.Range("M5").Value = 1.915838666666666
$ws.Range("N5").Value = 5.747515999999999
$ws.Range("O5").Value = 0.30709521276897
$ws.Range("P5").Value = 0.3070952127689701
$ws.Range("Q5").Value = 0.3359429488128889
$ws.Range("R5").Value = 3.023486539316
$ws.Range("S5").Value = 0.0006327465943822065
$ws.Range("T5").Value = 0.0006594822139217725

# Row 6
$ws.Range("I6").Value = 0.7783767285892832
$ws.Range("J6").Value = 0.8112656990851791
$ws.Range("M6").Value = 1.772429333333333
$ws.Range("N6").Value = 5.317288
$ws.Range("O6").Value = 0.2841077240522499
$ws.Range("P6").Value = 0.2841077240522499
$ws.Range("Q6").Value = 117.4109488881529
$ws.Range("R6").Value = 1056.698539993376
$ws.Range("S6").Value = 0.2211428408147371
$ws.Range("T6").Value = 0.2304868513687477

# Row 7
$ws.Range("I7").Value = 0.7783767285892832
$ws.Range("J7").Value = 0.8112656990851791
$ws.Range("O7").Value = 0.01265586741773478
$ws.Range("P7").Value = 0.01265586741773478
$ws.Range("S7").Value = 0.009851032678076094
$ws.Range("T7").Value = 0.01026727112817795

# Row 8
$ws.Range("I8").Value = 0.7783767285892832
$ws.Range("J8").Value = 0.8112656990851791
$ws.Range("M8").Value = 2.471359333333333
$ws.Range("N8").Value = 7.414078
$ws.Range("O8").Value = 0.3961411957610453
$ws.Range("P8").Value = 0.3961411957610453
$ws.Range("Q8").Value = 163.7101343976062
$ws.Range("R8").Value = 1473.391209578456
$ws.Range("S8").Value = 0.3083470880159292
$ws.Range("T8").Value = 0.3213757641155232

# Row 9
$ws.Range("I9").Value = 0.7783767285892832
$ws.Range("J9").Value = 0.8112656990851791
$ws.Range("M9").Value = 1.915838666666666
$ws.Range("N9").Value = 5.747515999999999
$ws.Range("O9").Value = 0.30709521276897
$ws.Range("P9").Value = 0.3070952127689701
$ws.Range("Q9").Value = 126.9108062812924
$ws.Range("R9").Value = 1142.197256531632
$ws.Range("S9").Value = 0.2390357670805408
$ws.Range("T9").Value = 0.2491358124727303

# Row 10
$ws.Range("G10").Value = 0.8839399999999999
$ws.Range("H10").Value = 2.65182
$ws.Range("I10").Value = 0.01038658934339153
$ws.Range("J10").Value = 0.01082545681966733
$ws.Range("M10").Value = 1.772429333333333
$ws.Range("N10").Value = 5.317288
$ws.Range("O10").Value = 0.2841077240522499
$ws.Range("P10").Value = 0.2841077240522499
$ws.Range("Q10").Value = 1.566721184906666
$ws.Range("R10").Value = 14.10049066416
$ws.Range("S10").Value = 0.002950910259016321
$ws.Range("T10").Value = 0.003075595898861593

# Row 11
$ws.Range("G11").Value = 0.8839399999999999
$ws.Range("H11").Value = 2.65182
$ws.Range("I11").Value = 0.01038658934339153
$ws.Range("J11").Value = 0.01082545681966733
$ws.Range("O11").Value = 0.01265586741773478
$ws.Range("P11").Value = 0.01265586741773478
$ws.Range("Q11").Value = 0.06979118805333333
$ws.Range("R11").Value = 0.6281206924799999
$ws.Range("S11").Value = 0.0001314512976524201
$ws.Range("T11").Value = 0.0001370055462461225

# Row 12
$ws.Range("G12").Value = 0.8839399999999999
$ws.Range("H12").Value = 2.65182
$ws.Range("I12").Value = 0.01038658934339153
$ws.Range("J12").Value = 0.01082545681966733
$ws.Range("M12").Value = 2.471359333333333
$ws.Range("N12").Value = 7.414078
$ws.Range("O12").Value = 0.3961411957610453
$ws.Range("P12").Value = 0.3961411957610453
$ws.Range("Q12").Value = 2.184533369106667
$ws.Range("R12").Value = 19.66080032196
$ws.Range("S12").Value = 0.004114555922370052
$ws.Range("T12").Value = 0.00428840940920258

# Row 13
$ws.Range("G13").Value = 0.8839399999999999
$ws.Range("H13").Value = 2.65182
$ws.Range("I13").Value = 0.01038658934339153
$ws.Range("J13").Value = 0.01082545681966733
$ws.Range("M13").Value = 1.915838666666666
$ws.Range("N13").Value = 5.747515999999999
$ws.Range("O13").Value = 0.30709521276897
$ws.Range("P13").Value = 0.3070952127689701
$ws.Range("Q13").Value = 1.693486431013333
$ws.Range("R13").Value = 15.24137787912
$ws.Range("S13").Value = 0.003189671864352739
$ws.Range("T13").Value = 0.003324445965357037

# Row 14
$ws.Range("G14").Value = 10.3504265
$ws.Range("H14").Value = 20.700853
$ws.Range("I14").Value = 0.1216209579659901
$ws.Range("J14").Value = 0.0845065616375851
$ws.Range("M14").Value = 1.772429333333333
$ws.Range("N14").Value = 5.317288
$ws.Range("O14").Value = 0.2841077240522499
$ws.Range("P14").Value = 0.2841077240522499
$ws.Range("Q14").Value = 18.34539954111066
$ws.Range("R14").Value = 110.072397246664
$ws.Range("S14").Value = 0.0345534535647718
$ws.Range("T14").Value = 0.02400896689433548

# Row 15
$ws.Range("G15").Value = 10.3504265
$ws.Range("H15").Value = 20.700853
$ws.Range("I15").Value = 0.1216209579659901
$ws.Range("J15").Value = 0.0845065616375851
$ws.Range("O15").Value = 0.01265586741773478
$ws.Range("P15").Value = 0.01265586741773478
$ws.Range("Q15").Value = 0.8172144741653332
$ws.Range("R15").Value = 4.903286844992
$ws.Range("S15").Value = 0.001539218719235465
$ws.Range("T15").Value = 0.001069503840013909

# Row 16
$ws.Range("G16").Value = 10.3504265
$ws.Range("H16").Value = 20.700853
$ws.Range("I16").Value = 0.1216209579659901
$ws.Range("J16").Value = 0.0845065616375851
$ws.Range("M16").Value = 2.471359333333333
$ws.Range("N16").Value = 7.414078
$ws.Range("O16").Value = 0.3961411957610453
$ws.Range("P16").Value = 0.3961411957610453
$ws.Range("Q16").Value = 25.57962313475567
$ws.Range("R16").Value = 153.477738808534
$ws.Range("S16").Value = 0.04817907171825115
$ws.Range("T16").Value = 0.03347653037676745

# Row 17
$ws.Range("G17").Value = 10.3504265
$ws.Range("H17").Value = 20.700853
$ws.Range("I17").Value = 0.1216209579659901
$ws.Range("J17").Value = 0.0845065616375851
$ws.Range("M17").Value = 1.915838666666666
$ws.Range("N17").Value = 5.747515999999999
$ws.Range("O17").Value = 0.30709521276897
$ws.Range("P17").Value = 0.3070952127689701
$ws.Range("Q17").Value = 19.82974730519133
$ws.Range("R17").Value = 118.978483831148
$ws.Range("S17").Value = 0.03734921396373169
$ws.Range("T17").Value = 0.02595156052646828

# Row 18
$ws.Range("G18").Value = 7.451303666666667
$ws.Range("H18").Value = 22.353911
$ws.Range("I18").Value = 0.08755529929471939
$ws.Range("J18").Value = 0.09125479794299257
$ws.Range("M18").Value = 1.772429333333333
$ws.Range("N18").Value = 5.317288
$ws.Range("O18").Value = 0.2841077240522499
$ws.Range("P18").Value = 0.2841077240522499
$ws.Range("Q18").Value = 13.20690919037422
$ws.Range("R18").Value = 118.862182713368
$ws.Range("S18").Value = 0.02487513681133629
$ws.Range("T18").Value = 0.02592619295243156

# Row 19
$ws.Range("G19").Value = 7.451303666666667
$ws.Range("H19").Value = 22.353911
$ws.Range("I19").Value = 0.08755529929471939
$ws.Range("J19").Value = 0.09125479794299257
$ws.Range("O19").Value = 0.01265586741773478
$ws.Range("P19").Value = 0.01265586741773478
$ws.Range("Q19").Value = 0.5883151972337778
$ws.Range("R19").Value = 5.294836775104
$ws.Range("S19").Value = 0.001108088259594056
$ws.Range("T19").Value = 0.00115490862399869

# Row 20
$ws.Range("G20").Value = 7.451303666666667
$ws.Range("H20").Value = 22.353911
$ws.Range("I20").Value = 0.08755529929471939
$ws.Range("J20").Value = 0.09125479794299257
$ws.Range("M20").Value = 2.471359333333333
$ws.Range("N20").Value = 7.414078
$ws.Range("O20").Value = 0.3961411957610453
$ws.Range("P20").Value = 0.3961411957610453
$ws.Range("Q20").Value = 18.41484886211756
$ws.Range("R20").Value = 165.733639759058
$ws.Range("S20").Value = 0.03468426095782635
$ws.Range("T20").Value = 0.03614978477606966

# Row 21
$ws.Range("G21").Value = 7.451303666666667
$ws.Range("H21").Value = 22.353911
$ws.Range("I21").Value = 0.08755529929471939
$ws.Range("J21").Value = 0.09125479794299257
$ws.Range("M21").Value = 1.915838666666666
$ws.Range("N21").Value = 5.747515999999999
$ws.Range("O21").Value = 0.30709521276897
$ws.Range("P21").Value = 0.3070952127689701
$ws.Range("Q21").Value = 14.27549568167511
$ws.Range("R21").Value = 128.479461135076
$ws.Range("S21").Value = 0.0268878132659627
$ws.Range("T21").Value = 0.02802391159049267
